$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("outcomes")

$ws.Range("A3").Value = 1089
$ws.Range("B3").Value = "Non-hemorrhagic Stroke"
$ws.Range("C3").Value = 365
